# 23rd commit - Core alerts - Placeholder and template management TCs
#
# The workbook originally had two tabs:
#   "Sheet2" - CollectionAgency/login related test-case rows
#   "Sheet1" - login test-case rows
# This edit removes the "Sheet2" tab entirely and repurposes the
# remaining "Sheet1" tab's rows 2-3 for the new Core Alerts &
# Notifications test classes, clearing out row 4's old content.

$wb = $excel.ActiveWorkbook

# Remove the "Sheet2" worksheet completely.
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null

$sheet1 = $wb.Worksheets.Item("Sheet1")

# Update the test-case names for rows 2 and 3 (Execute Flag stays "Yes").
$sheet1.Range("A2").Value = "Core.AlertsandNotifications.AlertsPlaceholderManagement_TestClass"
$sheet1.Range("A3").Value = "Core.AlertsandNotifications.AlertsTemplateManagement_TestClass"

# Row 4 no longer holds a test case - clear its contents (keep formatting).
$sheet1.Range("A4:B4").ClearContents() | Out-Null

# Make "Sheet1" the active/selected sheet with cell A9 selected, matching
# the saved view state of the workbook after the edit.
$sheet1.Activate() | Out-Null
$sheet1.Range("A9").Select() | Out-Null
